$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 181-182; this shifts the existing data (old rows
# 181..214) down to 183..216, matching the diff's row-shift pattern.
$ws.Range("A181:A182").EntireRow.Insert()

# Populate the new row 181 (weekly record: Primera quality)
$ws.Cells.Item(181, 1).Value = 4
$ws.Cells.Item(181, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(181, 3).Value = "Los Lagos"
$ws.Cells.Item(181, 4).Value = 44504
$ws.Cells.Item(181, 5).Value = 10
$ws.Cells.Item(181, 6).Value = 100112008
$ws.Cells.Item(181, 7).Value = "Coliflor"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 350
$ws.Cells.Item(181, 11).Value = 1200
$ws.Cells.Item(181, 12).Value = 1200
$ws.Cells.Item(181, 13).Value = 1200
$ws.Cells.Item(181, 14).Value = "$/unidad"
$ws.Cells.Item(181, 15).Value = "Región Metropolitana"
$ws.Cells.Item(181, 16).Value = 1200
$ws.Cells.Item(181, 17).Value = 1
$ws.Cells.Item(181, 18).Value = "Hortaliza"

# Populate the new row 182 (weekly record: Segunda quality, same date)
$ws.Cells.Item(182, 1).Value = 4
$ws.Cells.Item(182, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(182, 3).Value = "Los Lagos"
$ws.Cells.Item(182, 4).Value = 44504
$ws.Cells.Item(182, 5).Value = 10
$ws.Cells.Item(182, 6).Value = 100112008
$ws.Cells.Item(182, 7).Value = "Coliflor"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Segunda"
$ws.Cells.Item(182, 10).Value = 350
$ws.Cells.Item(182, 11).Value = 1000
$ws.Cells.Item(182, 12).Value = 1000
$ws.Cells.Item(182, 13).Value = 1000
$ws.Cells.Item(182, 14).Value = "$/unidad"
$ws.Cells.Item(182, 15).Value = "Región Metropolitana"
$ws.Cells.Item(182, 16).Value = 1000
$ws.Cells.Item(182, 17).Value = 1
$ws.Cells.Item(182, 18).Value = "Hortaliza"
